# Add two new student record rows (13 and 14) to the certificates sheet,
# mirroring the sparse layout already used by existing rows (e.g. row 12):
# only HALLTICKET (A), NAME (B) and STATUS (G) are populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: Demo Student / currently studying
$ws.Range("A13").Value = "HT2025EX"
$ws.Range("B13").Value = "Demo Student"
$ws.Range("G13").Value = "STUDYING"

# Row 14: Test2 / completed
$ws.Range("A14").Value = "TEST002"
$ws.Range("B14").Value = "Test2"
$ws.Range("G14").Value = "COMPLETED"
